$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (D and E) to make room for ownTeam/oppTeam;
# this shifts the existing batsman..sr columns from D:I to F:K.
$ws.Range("D1:E1").EntireColumn.Insert()

# Format the numeric-looking result columns (totalRuns..sr) as text so that
# values such as "81" or "158.82" are stored as strings, matching the source data.
$ws.Range("G1:K12").NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "venue"
$ws.Range("B1").Value = "date"
$ws.Range("C1").Value = "result"
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"
$ws.Range("F1").Value = "batsman"
$ws.Range("G1").Value = "totalRuns"
$ws.Range("H1").Value = "totalBalls"
$ws.Range("I1").Value = "total4s"
$ws.Range("J1").Value = "total6s"
$ws.Range("K1").Value = "sr"

# Row 2
$ws.Range("A2").Value = " Abu Dhabi"
$ws.Range("B2").Value = " October 07 2020"
$ws.Range("C2").Value = "KKR won by 10 runs"
$ws.Range("D2").Value = "Kolkata Knight Riders"
$ws.Range("E2").Value = "Chennai Super Kings"
$ws.Range("F2").Value = "Rahul Tripathi "
$ws.Range("G2").Value = "81"
$ws.Range("H2").Value = "51"
$ws.Range("I2").Value = "8"
$ws.Range("J2").Value = "3"
$ws.Range("K2").Value = "158.82"

# Row 3
$ws.Range("A3").Value = " Abu Dhabi"
$ws.Range("B3").Value = " October 16 2020"
$ws.Range("C3").Value = "Mumbai won by 8 wickets (with 19 balls remaining)"
$ws.Range("D3").Value = "Kolkata Knight Riders"
$ws.Range("E3").Value = "Mumbai Indians"
$ws.Range("F3").Value = "Rahul Tripathi "
$ws.Range("G3").Value = "7"
$ws.Range("H3").Value = "9"
$ws.Range("I3").Value = "1"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "77.77"

# Row 4
$ws.Range("A4").Value = " Abu Dhabi"
$ws.Range("B4").Value = " October 18 2020"
$ws.Range("C4").Value = "Match tied (KKR won the one-over eliminator)"
$ws.Range("D4").Value = "Kolkata Knight Riders"
$ws.Range("E4").Value = "Sunrisers Hyderabad"
$ws.Range("F4").Value = "Rahul Tripathi "
$ws.Range("G4").Value = "23"
$ws.Range("H4").Value = "16"
$ws.Range("I4").Value = "2"
$ws.Range("J4").Value = "1"
$ws.Range("K4").Value = "143.75"

# Row 5
$ws.Range("A5").Value = " Abu Dhabi"
$ws.Range("B5").Value = " October 10 2020"
$ws.Range("C5").Value = "KKR won by 2 runs"
$ws.Range("D5").Value = "Kolkata Knight Riders"
$ws.Range("E5").Value = "Kings XI Punjab"
$ws.Range("F5").Value = "Rahul Tripathi "
$ws.Range("G5").Value = "4"
$ws.Range("H5").Value = "10"
$ws.Range("I5").Value = "1"
$ws.Range("J5").Value = "0"
$ws.Range("K5").Value = "40.00"

# Row 6
$ws.Range("A6").Value = " Dubai (DSC)"
$ws.Range("B6").Value = " November 01 2020"
$ws.Range("C6").Value = "KKR won by 60 runs"
$ws.Range("D6").Value = "Kolkata Knight Riders"
$ws.Range("E6").Value = "Rajasthan Royals"
$ws.Range("F6").Value = "Rahul Tripathi "
$ws.Range("G6").Value = "39"
$ws.Range("H6").Value = "34"
$ws.Range("I6").Value = "4"
$ws.Range("J6").Value = "2"
$ws.Range("K6").Value = "114.70"

# Row 7
$ws.Range("A7").Value = " Sharjah"
$ws.Range("B7").Value = " October 03 2020"
$ws.Range("C7").Value = "Capitals won by 18 runs"
$ws.Range("D7").Value = "Kolkata Knight Riders"
$ws.Range("E7").Value = "Delhi Capitals"
$ws.Range("F7").Value = "Rahul Tripathi "
$ws.Range("G7").Value = "36"
$ws.Range("H7").Value = "16"
$ws.Range("I7").Value = "3"
$ws.Range("J7").Value = "3"
$ws.Range("K7").Value = "225.00"

# Row 8
$ws.Range("A8").Value = " Dubai (DSC)"
$ws.Range("B8").Value = " October 29 2020"
$ws.Range("C8").Value = "Super Kings won by 6 wickets"
$ws.Range("D8").Value = "Kolkata Knight Riders"
$ws.Range("E8").Value = "Chennai Super Kings"
$ws.Range("F8").Value = "Rahul Tripathi "
$ws.Range("G8").Value = "3"
$ws.Range("H8").Value = "2"
$ws.Range("I8").Value = "0"
$ws.Range("J8").Value = "0"
$ws.Range("K8").Value = "150.00"

# Row 9
$ws.Range("A9").Value = " Sharjah"
$ws.Range("B9").Value = " October 26 2020"
$ws.Range("C9").Value = "Kings XI won by 8 wickets (with 7 balls remaining)"
$ws.Range("D9").Value = "Kolkata Knight Riders"
$ws.Range("E9").Value = "Kings XI Punjab"
$ws.Range("F9").Value = "Rahul Tripathi "
$ws.Range("G9").Value = "7"
$ws.Range("H9").Value = "4"
$ws.Range("I9").Value = "0"
$ws.Range("J9").Value = "1"
$ws.Range("K9").Value = "175.00"

# Row 10
$ws.Range("A10").Value = " Abu Dhabi"
$ws.Range("B10").Value = " October 21 2020"
$ws.Range("C10").Value = "RCB won by 8 wickets (with 39 balls remaining)"
$ws.Range("D10").Value = "Kolkata Knight Riders"
$ws.Range("E10").Value = "Royal Challengers Bangalore"
$ws.Range("F10").Value = "Rahul Tripathi "
$ws.Range("G10").Value = "1"
$ws.Range("H10").Value = "5"
$ws.Range("I10").Value = "0"
$ws.Range("J10").Value = "0"
$ws.Range("K10").Value = "20.00"

# Row 11
$ws.Range("A11").Value = " Abu Dhabi"
$ws.Range("B11").Value = " October 24 2020"
$ws.Range("C11").Value = "KKR won by 59 runs"
$ws.Range("D11").Value = "Kolkata Knight Riders"
$ws.Range("E11").Value = "Delhi Capitals"
$ws.Range("F11").Value = "Rahul Tripathi "
$ws.Range("G11").Value = "13"
$ws.Range("H11").Value = "12"
$ws.Range("I11").Value = "1"
$ws.Range("J11").Value = "0"
$ws.Range("K11").Value = "108.33"

# Row 12
$ws.Range("A12").Value = " Sharjah"
$ws.Range("B12").Value = " October 12 2020"
$ws.Range("C12").Value = "RCB won by 82 runs"
$ws.Range("D12").Value = "Kolkata Knight Riders"
$ws.Range("E12").Value = "Royal Challengers Bangalore"
$ws.Range("F12").Value = "Rahul Tripathi "
$ws.Range("G12").Value = "16"
$ws.Range("H12").Value = "22"
$ws.Range("I12").Value = "1"
$ws.Range("J12").Value = "0"
$ws.Range("K12").Value = "72.72"

